# TalDoor_BOM.xlsx — mainboard PCB completed; P1/P2 generic 2x18 header
# replaced by a single bespoke P1 connector mated to a PocketBeagle.
# (commit: "Completed Mainboard PCB, needs Vgroove/Mousebite connection
#  between two. Also check motor driver")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: was the P1,P2 double-row 2x18 pin header; now a single P1
#     bespoke male header that plugs straight into a PocketBeagle. ---
# Clear the cells that no longer apply (Digikey part # / unit price) before
# writing the new values, so the row ends up with exactly the remaining
# columns populated (A,B,C,D,E,G,H,I — F empty).
$ws.Range("F18").ClearContents()
$ws.Range("G18").ClearContents()

$ws.Range("A18").Value = "P1"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "PocketBeagle"
$ws.Range("I18").Value = "Bespoke Male Pin Headers"
$ws.Range("D18").Value = "TalDoor_Footprints:PocketBeagle_TalDoor"
$ws.Range("E18").Value = "double row, odd1/even2 numbering"

# --- Selection / scroll position, matching the author's saved view ---
$ws.Range("C18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
